# Applies the "search by product name" refactor described in the commit:
# - Contas sheet: the generic "kirxhzees7" search term becomes the real
#   product name "Amidffffrd"; the stray search-results URL that used to
#   live in L2 (and its hyperlink) is removed because searching no longer
#   requires pasting an element/URL into the data mass.
# - BuscaLupa / BuscaHome sheets: the long marketing product names are
#   replaced by the short product names actually typed into the search
#   box, a new "chiclete" sample is added, and the now-unused result-page
#   hyperlinks/URLs are removed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Contas"
# ---------------------------------------------------------------------
$wsContas = $wb.Worksheets.Item("Contas")

$wsContas.Range("A2").Value = "Amidffffrd"

$wsContas.Range("L2").Hyperlinks.Delete()
$wsContas.Range("L2").Value = ""

# ---------------------------------------------------------------------
# Sheet "BuscaLupa"
# ---------------------------------------------------------------------
$wsBuscaLupa = $wb.Worksheets.Item("BuscaLupa")

$wsBuscaLupa.Range("A1").Value = "Headset"
$wsBuscaLupa.Range("B1").Value = "Tenis"
$wsBuscaLupa.Range("C1").Value = "chiclete"

$wsBuscaLupa.Range("A2").Value = "Headset H390"

$wsBuscaLupa.Range("A3").Hyperlinks.Delete()
$wsBuscaLupa.Range("A3").Value = ""

# Selecting A2 mirrors the saved cursor position in the edited file; restore
# "Contas" as the active tab afterwards since selecting a range on another
# sheet activates it.
[void]$wsBuscaLupa.Range("A2").Select()
[void]$wsContas.Activate()

# ---------------------------------------------------------------------
# Sheet "BuscaHome"
# ---------------------------------------------------------------------
$wsBuscaHome = $wb.Worksheets.Item("BuscaHome")

$wsBuscaHome.Range("A1").Value = "Tablet 608"
$wsBuscaHome.Range("B1").Value = "Headset H390"

$wsBuscaHome.Range("A2").Hyperlinks.Delete()
$wsBuscaHome.Range("A2").Value = ""

$wsBuscaHome.Range("B2").Hyperlinks.Delete()
$wsBuscaHome.Range("B2").Value = ""
